$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 669 ("「ラヤンはねむれない」...") entirely; Excel will automatically
# shift all subsequent rows up by one and adjust the used range/dimension.
$ws.Rows.Item(669).Delete()
